# Glossario update: translate new spell entries and refresh the Magias sheet.
# ("Traduzido da linha 1712 ate o fim do SkWd_new_us.u16 - Atualizado o
#  glossario com nomes de magias.")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Magias")

# Full, sorted English/Portuguese term list for the spells glossary (B:C),
# now extended with the newly translated entries.
$magiasData = @(
    @("Cancel", "Anulação"),
    @("Curse", "Maldição"),
    @("Dark ", "Escuridão"),
    @("Deflect", "Acurácia"),
    @("Dizzy", "Vertigem"),
    @("Erase", "Dissolução"),
    @("Extract", "Absorção"),
    @("Extracta", "Absorçãoe"),
    @("Extractus", "Absorçãous"),
    @("Flara", "Chamae"),
    @("Flare", "Chama"),
    @("Flarus", "Chamaus"),
    @("Ground", "Terra"),
    @("Grounda", "Terrae"),
    @("Groundus", "Terraus"),
    @("Heal", "Cura"),
    @("Heala", "Curae"),
    @("Healus", "Curaus"),
    @("KO", "nocaute"),
    @("Panic", "Pânico"),
    @("Paralysis", "Paralisia"),
    @("Petrify", "Petrificação"),
    @("Poison", "Sono"),
    @("Previve", "Auto-Revive"),
    @("Quick", "Rapidez"),
    @("Quicka", "Rapideze"),
    @("Quickus", "Rapidezus"),
    @("Reflect", "Reflexo"),
    @("Reflecta", "Reflexa"),
    @("Regenera", "Regeneraçãoe"),
    @("Regenerate", "Regeneração"),
    @("Resist", "Imune"),
    @("Resista", "Imuna"),
    @("Resistus", "Imunus"),
    @("Revive", "Revive"),
    @("Shadow", "Sombra"),
    @("Shadowa", "Sombrae"),
    @("Shadowus", "Sombraus"),
    @("Shell", "Carapaça"),
    @("Shella", "Carapaçae"),
    @("Shellus", "Carapaçaus"),
    @("Shield", "Escudo"),
    @("Shielda", "Escuda"),
    @("Shieldus", "Escudous"),
    @("Shina", "Luza"),
    @("Shine", "Luz"),
    @("Sleep", "Sono"),
    @("Slow", "Lentidão"),
    @("Slowa", "Lentidãoe"),
    @("Slowus", "Lentidãous"),
    @("status ailments", "enfermidades"),
    @("Trapfloor", "Armadilha"),
    @("Trapfloora", "Armadilhae"),
    @("Trapfloorus", "Armadilhaus"),
    @("Wall", "Muralha"),
    @("Walla", "Muralhae"),
    @("Wallus", "Muralhaus"),
    @("Water", "Água"),
    @("Watera", "Aguae"),
    @("Waterus", "Aguaus"),
    @("Wind", "Vento"),
    @("Winda", "Ventoe"),
    @("Windus", "Vetous"),
    @("Zephyr", "Vitalidade"),
    @("Zephyra", "Vitalidadea"),
    @("Zephyrus", "Vitalidadeus")
)

$rowCount = $magiasData.Count
$colCount = 2

$values = New-Object 'object[,]' $rowCount,$colCount
for ($i = 0; $i -lt $rowCount; $i++) {
    $values[$i,0] = $magiasData[$i][0]
    $values[$i,1] = $magiasData[$i][1]
}

$targetRange = $ws.Range($ws.Cells.Item(2,2), $ws.Cells.Item(1 + $rowCount, 3))
$targetRange.Value = $values

# Keep the glossary sorted A-Z by the English term, matching the sheet's
# existing sort state.
$sortRange = $ws.Range($ws.Cells.Item(2,2), $ws.Cells.Item(1 + $rowCount, 3))
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range($ws.Cells.Item(2,2), $ws.Cells.Item(1 + $rowCount, 2)))
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 0
$ws.Sort.Apply()

$ws.Activate()
$ws.Range("C12").Select()

Write-Output "Magias sheet updated: $rowCount rows"
